# Appends 5 new fixtures (rows 58-62) to the croatia/hnl 2023-2024 sheet,
# matching the "Atualizado por script em 31-10-2023 15:01" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 58 ---
$ws.Range("A57:V57").Copy()
$ws.Range("A58:V58").PasteSpecial(-4122)
$ws.Cells.Item(58,1).Value = 57
$ws.Cells.Item(58,2).Value = "croatia"
$ws.Cells.Item(58,3).Value = "hnl"
$ws.Cells.Item(58,4).Value = "2023-2024"
$ws.Cells.Item(58,5).Value = 45226.75
$ws.Cells.Item(58,6).Value = "Slaven Belupo"
$ws.Cells.Item(58,7).Value = 0
$ws.Cells.Item(58,8).Value = "Hajduk Split"
$ws.Cells.Item(58,9).Value = 1
$ws.Cells.Item(58,10).Value = 5.71
$ws.Cells.Item(58,11).Value = "22/10/2023 20:15"
$ws.Cells.Item(58,12).Value = 5.76
$ws.Cells.Item(58,13).Value = "27/10/2023 17:51"
$ws.Cells.Item(58,14).Value = 4.2
$ws.Cells.Item(58,15).Value = "22/10/2023 20:15"
$ws.Cells.Item(58,16).Value = 3.83
$ws.Cells.Item(58,17).Value = "27/10/2023 17:51"
$ws.Cells.Item(58,18).Value = 1.55
$ws.Cells.Item(58,19).Value = "22/10/2023 20:15"
$ws.Cells.Item(58,20).Value = 1.63
$ws.Cells.Item(58,21).Value = "27/10/2023 17:51"
$ws.Cells.Item(58,22).Value = "https://www.betexplorer.com/football/croatia/hnl/slaven-belupo-hajduk-split/fJ13HKZC/"

# --- row 59 ---
$ws.Range("A57:V57").Copy()
$ws.Range("A59:V59").PasteSpecial(-4122)
$ws.Cells.Item(59,1).Value = 58
$ws.Cells.Item(59,2).Value = "croatia"
$ws.Cells.Item(59,3).Value = "hnl"
$ws.Cells.Item(59,4).Value = "2023-2024"
$ws.Cells.Item(59,5).Value = 45227.66666666666
$ws.Cells.Item(59,6).Value = "Rudes"
$ws.Cells.Item(59,7).Value = 0
$ws.Cells.Item(59,8).Value = "Varazdin"
$ws.Cells.Item(59,9).Value = 0
$ws.Cells.Item(59,10).Value = 2.84
$ws.Cells.Item(59,11).Value = "21/10/2023 18:13"
$ws.Cells.Item(59,12).Value = 2.91
$ws.Cells.Item(59,13).Value = "28/10/2023 15:32"
$ws.Cells.Item(59,14).Value = 3.39
$ws.Cells.Item(59,15).Value = "21/10/2023 18:13"
$ws.Cells.Item(59,16).Value = 3.5
$ws.Cells.Item(59,17).Value = "28/10/2023 15:32"
$ws.Cells.Item(59,18).Value = 2.36
$ws.Cells.Item(59,19).Value = "21/10/2023 18:13"
$ws.Cells.Item(59,20).Value = 2.41
$ws.Cells.Item(59,21).Value = "28/10/2023 15:32"
$ws.Cells.Item(59,22).Value = "https://www.betexplorer.com/football/croatia/hnl/rudes-varazdin/CfCbI0K6/"

# --- row 60 ---
$ws.Range("A57:V57").Copy()
$ws.Range("A60:V60").PasteSpecial(-4122)
$ws.Cells.Item(60,1).Value = 59
$ws.Cells.Item(60,2).Value = "croatia"
$ws.Cells.Item(60,3).Value = "hnl"
$ws.Cells.Item(60,4).Value = "2023-2024"
$ws.Cells.Item(60,5).Value = 45227.75694444445
$ws.Cells.Item(60,6).Value = "Osijek"
$ws.Cells.Item(60,7).Value = 3
$ws.Cells.Item(60,8).Value = "Istra 1961"
$ws.Cells.Item(60,9).Value = 1
$ws.Cells.Item(60,10).Value = 1.56
$ws.Cells.Item(60,11).Value = "22/10/2023 20:15"
$ws.Cells.Item(60,12).Value = 1.52
$ws.Cells.Item(60,13).Value = "28/10/2023 16:23"
$ws.Cells.Item(60,14).Value = 4.1
$ws.Cells.Item(60,15).Value = "22/10/2023 20:15"
$ws.Cells.Item(60,16).Value = 4.21
$ws.Cells.Item(60,17).Value = "28/10/2023 18:00"
$ws.Cells.Item(60,18).Value = 5.72
$ws.Cells.Item(60,19).Value = "22/10/2023 20:15"
$ws.Cells.Item(60,20).Value = 6.55
$ws.Cells.Item(60,21).Value = "28/10/2023 18:00"
$ws.Cells.Item(60,22).Value = "https://www.betexplorer.com/football/croatia/hnl/osijek-istra-1961/8M57GvkJ/"

# --- row 61 ---
$ws.Range("A57:V57").Copy()
$ws.Range("A61:V61").PasteSpecial(-4122)
$ws.Cells.Item(61,1).Value = 60
$ws.Cells.Item(61,2).Value = "croatia"
$ws.Cells.Item(61,3).Value = "hnl"
$ws.Cells.Item(61,4).Value = "2023-2024"
$ws.Cells.Item(61,5).Value = 45228.625
$ws.Cells.Item(61,6).Value = "Rijeka"
$ws.Cells.Item(61,7).Value = 1
$ws.Cells.Item(61,8).Value = "Gorica"
$ws.Cells.Item(61,9).Value = 0
$ws.Cells.Item(61,10).Value = 1.54
$ws.Cells.Item(61,11).Value = "22/10/2023 16:12"
$ws.Cells.Item(61,12).Value = 1.53
$ws.Cells.Item(61,13).Value = "29/10/2023 14:56"
$ws.Cells.Item(61,14).Value = 4.16
$ws.Cells.Item(61,15).Value = "22/10/2023 16:12"
$ws.Cells.Item(61,16).Value = 4.04
$ws.Cells.Item(61,17).Value = "29/10/2023 14:56"
$ws.Cells.Item(61,18).Value = 5.84
$ws.Cells.Item(61,19).Value = "22/10/2023 16:12"
$ws.Cells.Item(61,20).Value = 6.62
$ws.Cells.Item(61,21).Value = "29/10/2023 14:56"
$ws.Cells.Item(61,22).Value = "https://www.betexplorer.com/football/croatia/hnl/rijeka-hnk-gorica/MPgkuORJ/"

# --- row 62 ---
$ws.Range("A57:V57").Copy()
$ws.Range("A62:V62").PasteSpecial(-4122)
$ws.Cells.Item(62,1).Value = 61
$ws.Cells.Item(62,2).Value = "croatia"
$ws.Cells.Item(62,3).Value = "hnl"
$ws.Cells.Item(62,4).Value = "2023-2024"
$ws.Cells.Item(62,5).Value = 45228.71527777778
$ws.Cells.Item(62,6).Value = "D. Zagreb"
$ws.Cells.Item(62,7).Value = 2
$ws.Cells.Item(62,8).Value = "Lok. Zagreb"
$ws.Cells.Item(62,9).Value = 1
$ws.Cells.Item(62,10).Value = 1.25
$ws.Cells.Item(62,11).Value = "22/10/2023 20:15"
$ws.Cells.Item(62,12).Value = 1.42
$ws.Cells.Item(62,13).Value = "29/10/2023 17:02"
$ws.Cells.Item(62,14).Value = 5.87
$ws.Cells.Item(62,15).Value = "22/10/2023 20:15"
$ws.Cells.Item(62,16).Value = 4.69
$ws.Cells.Item(62,17).Value = "29/10/2023 17:02"
$ws.Cells.Item(62,18).Value = 10.65
$ws.Cells.Item(62,19).Value = "22/10/2023 20:15"
$ws.Cells.Item(62,20).Value = 7.55
$ws.Cells.Item(62,21).Value = "29/10/2023 17:02"
$ws.Cells.Item(62,22).Value = "https://www.betexplorer.com/football/croatia/hnl/din-zagreb-lok-zagreb/0GhgvrsQ/"

$excel.CutCopyMode = 0
